$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classrooms")

# Insert a new row before row 5 (shifts existing rows 5-16 down to 6-17)
$ws.Rows.Item(5).Insert()

# The fresh row picks up the column's default style (132 / "Normal 2") on
# insert; reset it back to Normal and to the sheet's standard row height so
# it matches the plain, unstyled look of the surrounding data rows.
$ws.Range("A5:B5").Style = "Normal"
$ws.Rows.Item(5).RowHeight = 15

# Populate the newly inserted row 5 with the new classroom entry
$ws.Cells.Item(5, 1).Value = "11-536"
$ws.Cells.Item(5, 2).Value = 40

# Update the selection to match the post-edit workbook state
$ws.Activate()
$ws.Range("D8").Select()
